$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-ASCII helper characters (built via code point so the file itself can
# stay plain ASCII and survive any tooling round-trip unmangled).
$nbsp   = [char]0x00A0
$eacute = [char]0x00E9

# ---------------------------------------------------------------------------
# Append 16 new training-log rows (367-382) for 2025-09-17 (serial 45917),
# mirroring the existing layout: A=Date, B=Player, C=Volume, D=Intensite,
# E=Charge, F=Douleur, G=Localisation douleur, H=Plaisir, I=C*D (formula).
# ---------------------------------------------------------------------------

# 1) Seed formatting for the whole new block from a template row that has
#    an empty "G" (localisation douleur) cell, so styles (s=3/1/2) match
#    without minting new style indices.
$ws.Range("A322:I322").Copy()
$ws.Range("A367:I382").PasteSpecial(-4122)

# 2) Rows whose "G" column is populated need the "populated" G style (s=1)
#    instead of the "empty" style (s=2); copy that from an existing
#    populated G cell.
$ws.Range("G366").Copy()
$ws.Range("G368").PasteSpecial(-4122)
$ws.Range("G369").PasteSpecial(-4122)
$ws.Range("G370").PasteSpecial(-4122)
$ws.Range("G375").PasteSpecial(-4122)
$ws.Range("G376").PasteSpecial(-4122)
$ws.Range("G377").PasteSpecial(-4122)
$ws.Range("G380").PasteSpecial(-4122)
$ws.Range("G381").PasteSpecial(-4122)

function Add-WellnessRow {
    param(
        [int]$Row,
        [string]$Player,
        [double]$Volume,
        [double]$Intensite,
        [double]$Charge,
        [double]$Douleur,
        $Localisation,
        [double]$Plaisir
    )

    $ws.Cells.Item($Row, 1).Value = 45917
    $ws.Cells.Item($Row, 2).Value = $Player
    $ws.Cells.Item($Row, 3).Value = $Volume
    $ws.Cells.Item($Row, 4).Value = $Intensite
    $ws.Cells.Item($Row, 5).Value = $Charge
    $ws.Cells.Item($Row, 6).Value = $Douleur
    if ($Localisation) {
        $ws.Cells.Item($Row, 7).Value = $Localisation
    }
    $ws.Cells.Item($Row, 8).Value = $Plaisir
    $ws.Cells.Item($Row, 9).Formula = "=C$Row*D$Row"
}

$locCheville369  = "Cheville" + $nbsp
$locAdducteur370 = "Adducteur" + $nbsp
$playerKarahali  = "Karahali Souar" + $eacute

Add-WellnessRow 367 "Sofiane Belle"   75 5 4  0 $null                5
Add-WellnessRow 368 "Amir Etien"      75 6 7  4 "Flanc"              7
Add-WellnessRow 369 "Yoan Zouma"      75 3 5  4 $locCheville369      5
Add-WellnessRow 370 "Ilyes Boughanmi" 75 5 5  4 $locAdducteur370     0
Add-WellnessRow 371 "Omar Benyounes"  75 5 5  0 $null                8
Add-WellnessRow 372 "Yanis Berrached" 75 6 6  0 $null                2
Add-WellnessRow 373 "Malik Boussaid"  75 2 0  0 $null                10
Add-WellnessRow 374 "Jeremie Laurent" 75 7 6  0 $null                7
Add-WellnessRow 375 "Kamal Bafounta"  75 5 0  1 "Genou cheville"     8
Add-WellnessRow 376 "Naim Ighbane"    75 6 6  2 "Cheville"           6
Add-WellnessRow 377 "Yoann Martelat"  75 4 4  8 "Genou"              4
Add-WellnessRow 378 "Ilan Ihaddadene" 75 4 4  0 $null                7
Add-WellnessRow 379 "Emmanuel Valey"  75 7 6  0 $null                1
Add-WellnessRow 380 $playerKarahali   75 2 10 6 "Cheville"           1
Add-WellnessRow 381 "Naim Dhib"       75 5 5  3 "Genou"              5
Add-WellnessRow 382 "Mattheo Haon"    75 5 5  0 $null                8

# ---------------------------------------------------------------------------
# Update the saved view state to match the authored edit (selection parked
# on the last-entered row).
# ---------------------------------------------------------------------------
$ws.Range("L376").Select() | Out-Null

$excel.Calculate()
